# Update column F (dSF) values for the specific rows per the commit:
# "repull data, push all data, mean calculation"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    8  = -1
    17 = -3
    22 = -1
    25 = 4
    27 = -10
    31 = -2
    37 = 2
    42 = -1
    43 = 4
    46 = -3
    53 = 2
    54 = -1
    57 = -1
    61 = -2
    62 = -2
    64 = -4
    65 = -4
    66 = 2
    67 = -1
    71 = 14
    75 = -5
    76 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
